$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "DID Read" (sheet1) — columns were C=Data/D=Error but values were
# actually being written into the wrong column whenever the Data payload was
# missing; fix the header order and move the misplaced values to column D,
# and mirror the write-DID fix for the "KO"/error case (NOK now carries the
# error text in column C).
# ---------------------------------------------------------------------------
$wsRead = $wb.Worksheets.Item("DID Read")

# Header row: swap "Data"/"Error" so column C = Error, column D = Data.
$wsRead.Range("C1").Value = "Error"
$wsRead.Range("D1").Value = "Data"

# Row 2 (DID 01E3): data payload was wrongly written to C2 -- move it to D2.
$wsRead.Range("D2").Value = "0x00;0x00;0x00;0x00;0x00;0x00;0x00;0x00;0x00;0x00"
$wsRead.Range("C2").ClearContents()

# Row 3 (DID 01E4): same fix -- move the 5-byte-repeated payload to D3.
$wsRead.Range("D3").Value = "0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00;0x80;0x00;0x05;0x00"
$wsRead.Range("C3").ClearContents()

# Row 4 (DID 8012): result "KO" -> "NOK", and the error text moves from D4 to C4.
$wsRead.Range("B4").Value = "NOK"
$wsRead.Range("C4").Value = "Negative response: Error code 0x14: Response too long"
$wsRead.Range("D4").ClearContents()

# Column width adaptation.
$wsRead.Columns.Item(1).ColumnWidth = 6 - 5/6
$wsRead.Columns.Item(2).ColumnWidth = 10 - 5/6
$wsRead.Columns.Item(3).ColumnWidth = 55 - 5/6
$wsRead.Columns.Item(4).ColumnWidth = 201 - 5/6

# Colour rule: green for OK, red for NOK on the Resultat column.
$rngReadStatus = $wsRead.Range("B2:B4")
$cfReadOk = $rngReadStatus.FormatConditions.Add(1, 3, '"OK"')
$cfReadOk.Interior.Color = 65280
$cfReadNok = $rngReadStatus.FormatConditions.Add(1, 3, '"NOK"')
$cfReadNok.Interior.Color = 255

# ---------------------------------------------------------------------------
# Sheet "DID Write" (sheet2) — same fix: rotate the Data/Status/Error headers
# back into the correct columns, and since a 5-byte write response wrote its
# status into the wrong columns, move each row's values over by one column.
# ---------------------------------------------------------------------------
$wsWrite = $wb.Worksheets.Item("DID Write")

# Header row: B=Status, C=Error, D=Data.
$wsWrite.Range("B1").Value = "Status"
$wsWrite.Range("C1").Value = "Error"
$wsWrite.Range("D1").Value = "Data"

# Row 2 (DID 8283): Status "Failed" -> "NOK", Error text stays conceptually
# but moves to C2, and the raw data "0;0" moves from B2 to D2.
$wsWrite.Range("B2").Value = "NOK"
$wsWrite.Range("C2").Value = "Negative response: Error code 0x13: Invalid message length/format"
$wsWrite.Range("D2").Value = "0;0"

# Row 3 (DID 043A): same fix, and the error text is corrected to match the
# real negative response (it used to incorrectly say "invalid index").
$wsWrite.Range("B3").Value = "NOK"
$wsWrite.Range("C3").Value = "Negative response: Error code 0x13: Invalid message length/format"
$wsWrite.Range("D3").Value = "0;0;0;0;1"

# Column width adaptation.
$wsWrite.Columns.Item(1).ColumnWidth = 6 - 5/6
$wsWrite.Columns.Item(2).ColumnWidth = 8 - 5/6
$wsWrite.Columns.Item(3).ColumnWidth = 67 - 5/6
$wsWrite.Columns.Item(4).ColumnWidth = 11 - 5/6

# Colour rule: green for OK, red for NOK on the Status column.
$rngWriteStatus = $wsWrite.Range("B2:B3")
$cfWriteOk = $rngWriteStatus.FormatConditions.Add(1, 3, '"OK"')
$cfWriteOk.Interior.Color = 65280
$cfWriteNok = $rngWriteStatus.FormatConditions.Add(1, 3, '"NOK"')
$cfWriteNok.Interior.Color = 255
